$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet - Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-17 01:00:40"

# zh-cn sheet - Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-17 01:00:36"
$wsZhCn.Range("K2").Value = "2016-08-17 01:00:53"

# de-de sheet - Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-08-17 01:00:40"
$wsDeDe.Range("K2").Value = "2016-08-17 01:01:02"
